$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.416.65"
$ws.Range("E2").Value = "  +0.64%  "
$ws.Range("D3").Value = "2.377.29"
$ws.Range("E3").Value = "  +0.42%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "553.07"
$ws.Range("E5").Value = "  +1.97%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.77"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  -0.24%  "
$ws.Range("D9").Value = "2.378.27"
$ws.Range("E9").Value = "  +0.47%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.108"
$ws.Range("E10").Value = "  +2.58%  "
$ws.Range("E11").Value = "  +2.12%  "
$ws.Range("E12").Value = "  +2.24%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.351"
$ws.Range("E13").Value = "  +2.59%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.50"
$ws.Range("E14").Value = "  +2.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000171"
$ws.Range("E15").Value = "  +5.39%  "
$ws.Range("D16").Value = "2.805.73"
$ws.Range("E16").Value = "  +0.53%  "
$ws.Range("D17").Value = "61.294.53"
$ws.Range("E17").Value = "  +0.71%  "
$ws.Range("D18").Value = "2.375.26"
$ws.Range("E18").Value = "  +0.59%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.98"
$ws.Range("E19").Value = "  +3.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.16"
$ws.Range("E20").Value = "  +2.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "320.80"
$ws.Range("E21").Value = "  +1.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.69"
$ws.Range("E22").Value = "  +1.34%  "
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("E24").Value = "  -8.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "64.22"
$ws.Range("E25").Value = "  +1.25%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.84"
$ws.Range("E26").Value = "  +3.99%  "
$ws.Range("E27").Value = "  -0.13%  "
$ws.Range("D28").Value = "2.492.71"
$ws.Range("E28").Value = "  +0.36%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.16"
$ws.Range("E29").Value = "  +2.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "519.42"
$ws.Range("E30").Value = "  +2.37%  "
$ws.Range("D31").Value = "0.0₃0903"
$ws.Range("E31").Value = "  -0.38%  "
$ws.Range("E32").Value = "  -0.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.148"
$ws.Range("E33").Value = "  +1.74%  "
$ws.Range("E34").Value = "  +2.78%  "
$ws.Range("E35").Value = "  -1.62%  "
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.53"
$ws.Range("E37").Value = "  +4.73%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.69"
$ws.Range("E38").Value = "  +2.14%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.89"
$ws.Range("E39").Value = "  +6.09%  "
$ws.Range("E40").Value = "  +1.53%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.49"
$ws.Range("E41").Value = "  -0.47%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "147.01"
$ws.Range("E42").Value = "  +6.16%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.37"
$ws.Range("E44").Value = "  +2.93%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "148.18"
$ws.Range("E45").Value = "  +6.62%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.15"
$ws.Range("E46").Value = "  +0.42%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.60"
$ws.Range("E47").Value = "  +2.06%  "
$ws.Range("E48").Value = "  +1.71%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.72"
$ws.Range("E49").Value = "  +0.63%  "
$ws.Range("E50").Value = "  +1.40%  "
$ws.Range("E51").Value = "  +1.02%  "
